# New crime data collected — weekly CompStat 122nd Precinct refresh.
# Moves the report window forward one week (7/8-7/14/2024 -> 7/15-7/21/2024,
# "Number 28" -> "Number 29") and replaces the crime-count table (rows 15-33)
# with the new week's figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Header text: volume/number banner and the "week covering" date range.
# ---------------------------------------------------------------------
$ws.Range("A8").Value = "Volume 31   Number  29"
$ws.Range("C9").Value = "Report Covering the Week  7/15/2024  Through  7/21/2024"

# ---------------------------------------------------------------------
# Helper: paste only the number-format/style of $srcA1 onto $dstA1
# without touching the value that was just written there. Mirrors what
# Excel does when a cell's content type changes (text <-> number) but a
# specific existing style needs to be (re)applied.
# ---------------------------------------------------------------------
function Copy-Style([string]$srcA1, [string]$dstA1) {
    $ws.Range($srcA1).Copy() | Out-Null
    $ws.Range($dstA1).PasteSpecial(-4122) | Out-Null
}

# ---------------------------------------------------------------------
# Row 15 - Rape
# ---------------------------------------------------------------------
$ws.Range("C15").Value = 1
Copy-Style "D16" "C15"
$ws.Range("F15").Value = 1
Copy-Style "G15" "F15"
$ws.Range("H15").Value = 0
$ws.Range("I15").Value = 10
$ws.Range("K15").Value = 66.666666666666
$ws.Range("L15").Value = 150
$ws.Range("M15").Value = 0
$ws.Range("N15").Value = -33.333333333333

# ---------------------------------------------------------------------
# Row 16 - Robbery
# ---------------------------------------------------------------------
$ws.Range("D16").Value = 3
$ws.Range("E16").Value = -66.666666666666
$ws.Range("G16").Value = 4
$ws.Range("H16").Value = 50
$ws.Range("I16").Value = 34
$ws.Range("J16").Value = 31
$ws.Range("K16").Value = 9.677419354838
$ws.Range("L16").Value = -2.857142857142
$ws.Range("M16").Value = -41.379310344827
$ws.Range("N16").Value = -78.343949044586

# ---------------------------------------------------------------------
# Row 17 - Fel. Assault
# ---------------------------------------------------------------------
$ws.Range("C17").Value = 2
$ws.Range("E17").Value = -33.333333333333
$ws.Range("F17").Value = 10
$ws.Range("G17").Value = 16
$ws.Range("H17").Value = -37.5
$ws.Range("I17").Value = 93
$ws.Range("J17").Value = 105
$ws.Range("K17").Value = -11.428571428571
$ws.Range("L17").Value = 66.071428571428
$ws.Range("M17").Value = 20.779220779220
$ws.Range("N17").Value = -41.875

# ---------------------------------------------------------------------
# Row 18 - Burglary
# ---------------------------------------------------------------------
$ws.Range("C18").Value = 1
Copy-Style "D18" "C18"
$ws.Range("D18").Value = 6
$ws.Range("E18").Value = -83.333333333333
$ws.Range("F18").Value = 3
$ws.Range("G18").Value = 11
$ws.Range("H18").Value = -72.727272727272
$ws.Range("I18").Value = 46
$ws.Range("J18").Value = 58
$ws.Range("K18").Value = -20.689655172413
$ws.Range("L18").Value = 21.052631578947
$ws.Range("M18").Value = -57.798165137614
$ws.Range("N18").Value = -93.866666666666

# ---------------------------------------------------------------------
# Row 19 - Gr. Larceny
# ---------------------------------------------------------------------
$ws.Range("C19").Value = 6
$ws.Range("D19").Value = 16
$ws.Range("E19").Value = -62.5
$ws.Range("F19").Value = 33
$ws.Range("G19").Value = 36
$ws.Range("H19").Value = -8.333333333333
$ws.Range("I19").Value = 282
$ws.Range("J19").Value = 267
$ws.Range("K19").Value = 5.617977528089
$ws.Range("L19").Value = 55.801104972375
$ws.Range("M19").Value = 25.333333333333
$ws.Range("N19").Value = -38.695652173913

# ---------------------------------------------------------------------
# Row 20 - G.L.A.
# ---------------------------------------------------------------------
$ws.Range("C20").Value = 4
$ws.Range("E20").Value = 0
$ws.Range("F20").Value = 15
$ws.Range("G20").Value = 12
$ws.Range("H20").Value = 25
$ws.Range("I20").Value = 47
$ws.Range("J20").Value = 58
$ws.Range("K20").Value = -18.965517241379
$ws.Range("L20").Value = -32.857142857142
$ws.Range("M20").Value = -14.545454545454
$ws.Range("N20").Value = -96.926095487246

# ---------------------------------------------------------------------
# Row 21 - TOTAL
# ---------------------------------------------------------------------
$ws.Range("C21").Value = 15
$ws.Range("D21").Value = 32
$ws.Range("E21").Value = -53.125
$ws.Range("F21").Value = 68
$ws.Range("G21").Value = 80
$ws.Range("H21").Value = -15
$ws.Range("I21").Value = 512
$ws.Range("J21").Value = 526
$ws.Range("K21").Value = -2.661596958174
$ws.Range("L21").Value = 32.987012987013
$ws.Range("M21").Value = -4.119850187265
$ws.Range("N21").Value = -83.338756915066

# ---------------------------------------------------------------------
# Row 23 - Housing
# ---------------------------------------------------------------------
$ws.Range("C23").Value = 1
Copy-Style "F23" "C23"
$ws.Range("I23").Value = 11
$ws.Range("J23").Value = 21
$ws.Range("K23").Value = -47.619047619047
$ws.Range("L23").Value = -8.333333333333
$ws.Range("M23").Value = 10

# ---------------------------------------------------------------------
# Row 24 - Petit Larceny
# ---------------------------------------------------------------------
$ws.Range("C24").Value = 17
$ws.Range("D24").Value = 20
$ws.Range("E24").Value = -15
$ws.Range("F24").Value = 81
$ws.Range("H24").Value = -6.896551724137
$ws.Range("I24").Value = 576
$ws.Range("J24").Value = 621
$ws.Range("K24").Value = -7.246376811594
$ws.Range("L24").Value = 44.723618090452
$ws.Range("M24").Value = -39.368421052631

# ---------------------------------------------------------------------
# Row 25 - Retail Theft
# ---------------------------------------------------------------------
$ws.Range("C25").Value = 4
$ws.Range("D25").Value = 12
$ws.Range("E25").Value = -66.666666666666
$ws.Range("F25").Value = 29
$ws.Range("G25").Value = 28
$ws.Range("H25").Value = 3.571428571428
$ws.Range("I25").Value = 269
$ws.Range("J25").Value = 266
$ws.Range("K25").Value = 1.127819548872
$ws.Range("L25").Value = 192.391304347826

# ---------------------------------------------------------------------
# Row 26 - Misd. Assault
# ---------------------------------------------------------------------
$ws.Range("C26").Value = 7
$ws.Range("D26").Value = 5
$ws.Range("E26").Value = 40
$ws.Range("F26").Value = 32
$ws.Range("H26").Value = 28
$ws.Range("I26").Value = 204
$ws.Range("J26").Value = 173
$ws.Range("K26").Value = 17.919075144508
$ws.Range("L26").Value = 4.081632653061
$ws.Range("M26").Value = -37.037037037037

# ---------------------------------------------------------------------
# Row 27 - UCR Rape*
# ---------------------------------------------------------------------
$ws.Range("C27").Value = 1
Copy-Style "G27" "C27"
$ws.Range("F27").Value = 1
Copy-Style "G27" "F27"
$ws.Range("H27").Value = 0
$ws.Range("I27").Value = 13
$ws.Range("K27").Value = 44.444444444444
$ws.Range("L27").Value = 18.181818181818

# ---------------------------------------------------------------------
# Row 28 - Other Sex Crimes (some numeric cells become "no data" text)
# ---------------------------------------------------------------------
$ws.Range("C28").Value = "'0"
Copy-Style "C29" "C28"
$ws.Range("D28").Value = "'0"
Copy-Style "D29" "D28"
$ws.Range("E28").Value = "***.*"
Copy-Style "M28" "E28"
$ws.Range("F28").Value = 2
$ws.Range("G28").Value = 3
$ws.Range("I28").Value = 22
$ws.Range("K28").Value = 0
$ws.Range("L28").Value = 10

# ---------------------------------------------------------------------
# Row 29 - Shooting Vic. (28-day count becomes "no data")
# ---------------------------------------------------------------------
$ws.Range("F29").Value = "'0"
Copy-Style "G29" "F29"

# ---------------------------------------------------------------------
# Row 30 - Shooting Inc. (28-day count becomes "no data")
# ---------------------------------------------------------------------
$ws.Range("F30").Value = "'0"
Copy-Style "G30" "F30"

# ---------------------------------------------------------------------
# Row 31 - Hate Crimes (28-day count becomes "no data"; 2-year % updates)
# ---------------------------------------------------------------------
$ws.Range("F31").Value = "'0"
Copy-Style "G31" "F31"
$ws.Range("L31").Value = -33.333333333333

# ---------------------------------------------------------------------
# Row 33 - Traffic Fatalities (2-year % now has a value instead of "no data")
# ---------------------------------------------------------------------
$ws.Range("L33").Value = 0
Copy-Style "K33" "L33"
